$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "编号" (ID/number) column A entirely, shifting B,C,D,E left
# into A,B,C,D. This matches the new layout: 项目名称 / 项目代码 / 提示 / 结果
$xlShiftToLeft = -4159
$ws.Range("A1:A14").Delete($xlShiftToLeft)

# Clean up the leading numeric index that used to be concatenated onto the
# "项目名称" text (now living in column A) for several rows, and correct a
# couple of mislabeled project names.
$ws.Range("A3").Value = "尿蛋白"
$ws.Range("A4").Value = "胆红素"
$ws.Range("A7").Value = "比重"
$ws.Range("A9").Value = "耐体"
$ws.Range("A10").Value = "亚硝酸盐"
$ws.Range("A11").Value = "白细胞"
$ws.Range("A12").Value = "RDW-CV"
$ws.Range("A13").Value = "外观"
$ws.Range("A14").Value = "红细胞(镜检)"
